{"js": "const replacements = [\n  [\"645\u00f79=\", \"913\u00f72=\"],\n  [\"809\u00f73=\", \"137\u00f78=\"],\n  [\"218\u00f73=\", \"993\u00f75=\"],\n  [\"170\u00f73=\", \"343\u00f77=\"],\n  [\"313\u00f73=\", \"252\u00f76=\"],\n  [\"378\u00f79=\", \"297\u00f78=\"],\n  [\"921\u00f77=\", \"358\u00f76=\"],\n  [\"362\u00f77=\", \"509\u00f72=\"],\n  [\"640\u00f76=\", \"131\u00f78=\"],\n  [\"929\u00f74=\", \"964\u00f75=\"],\n  [\"805\u00f72=\", \"463\u00f72=\"],\n  [\"632\u00f79=\", \"505\u00f73=\"],\n  [\"122\u00f78=\", \"110\u00f73=\"],\n  [\"556\u00f77=\", \"769\u00f79=\"],\n  [\"395\u00f73=\", \"748\u00f76=\"],\n  [\"789\u00f76=\", \"666\u00f74=\"],\n  [\"396\u00f73=\", \"297\u00f76=\"],\n  [\"132\u00f72=\", \"950\u00f75=\"],\n  [\"843\u00f78=\", \"627\u00f77=\"],\n  [\"584\u00f74=\", \"179\u00f73=\"],\n  [\"167\u00f72=\", \"160\u00f73=\"],\n  [\"633\u00f74=\", \"856\u00f78=\"],\n  [\"773\u00f74=\", \"814\u00f75=\"],\n  [\"537\u00f72=\", \"934\u00f77=\"],\n  [\"266\u00f75=\", \"439\u00f77=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"645\u00f79=\", \"913\u00f72=\"),\n    @(\"809\u00f73=\", \"137\u00f78=\"),\n    @(\"218\u00f73=\", \"993\u00f75=\"),\n    @(\"170\u00f73=\", \"343\u00f77=\"),\n    @(\"313\u00f73=\", \"252\u00f76=\"),\n    @(\"378\u00f79=\", \"297\u00f78=\"),\n    @(\"921\u00f77=\", \"358\u00f76=\"),\n    @(\"362\u00f77=\", \"509\u00f72=\"),\n    @(\"640\u00f76=\", \"131\u00f78=\"),\n    @(\"929\u00f74=\", \"964\u00f75=\"),\n    @(\"805\u00f72=\", \"463\u00f72=\"),\n    @(\"632\u00f79=\", \"505\u00f73=\"),\n    @(\"122\u00f78=\", \"110\u00f73=\"),\n    @(\"556\u00f77=\", \"769\u00f79=\"),\n    @(\"395\u00f73=\", \"748\u00f76=\"),\n    @(\"789\u00f76=\", \"666\u00f74=\"),\n    @(\"396\u00f73=\", \"297\u00f76=\"),\n    @(\"132\u00f72=\", \"950\u00f75=\"),\n    @(\"843\u00f78=\", \"627\u00f77=\"),\n    @(\"584\u00f74=\", \"179\u00f73=\"),\n    @(\"167\u00f72=\", \"160\u00f73=\"),\n    @(\"633\u00f74=\", \"856\u00f78=\"),\n    @(\"773\u00f74=\", \"814\u00f75=\"),\n    @(\"537\u00f72=\", \"934\u00f77=\"),\n    @(\"266\u00f75=\", \"439\u00f77=\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
